$d = $word.ActiveDocument

# Locate the paragraph that contains "Ver no Jupiter ..." - the footer block to
# be removed consists of three paragraphs:
#   1) an empty paragraph right before it,
#   2) "Ver no Jupiter Salvar em pdf Salvar em docx"
#   3) the "... Contact: luizeleno@usp.br ..." copyright paragraph right after it
# All three must be removed as a single unit, leaving the surrounding paragraphs
# (the "LOQ4064: ..." requisitos line before, and the blank paragraph that
# precedes the page-break paragraph after) untouched.

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text
    if ($text -match "Ver no Jupiter Salvar em pdf Salvar em docx") {
        $blankBefore = $d.Paragraphs.Item($i - 1)
        $copyrightAfter = $d.Paragraphs.Item($i + 1)

        $startPos = $blankBefore.Range.Start
        $endPos = $copyrightAfter.Range.End

        $killRange = $d.Range($startPos, $endPos)
        $killRange.Delete()
        break
    }
}
